# [Feat 2269] Add support of parameters worksheet metadata.
# Rename "PARAMETERS TODO" -> "PARAMETERS", populate its header row
# (mirroring the STEPS header style) and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Rename the "PARAMETERS TODO" sheet to "PARAMETERS"
$paramSheet = $wb.Worksheets.Item("PARAMETERS TODO")
$paramSheet.Name = "PARAMETERS"

# Populate the header row with the parameter metadata columns
$paramSheet.Range("A1").Value = "ACTION"
$paramSheet.Range("B1").Value = "TC_OWNER_PATH"
$paramSheet.Range("C1").Value = "TC_OWNER_ID"
$paramSheet.Range("D1").Value = "TC_PARAM_ID"
$paramSheet.Range("E1").Value = "TC_PARAM_NAME"
$paramSheet.Range("F1").Value = "TC_PARAM_DESCRIPTION"

# Auto-fit the header columns so their width matches the other metadata sheets
# (widths mirror the "best fit" values Excel computes for these header labels)
$paramSheet.Range("B1").ColumnWidth = 16.022135416666668
$paramSheet.Range("C1").ColumnWidth = 13.022135416666666
$paramSheet.Range("D1").ColumnWidth = 12.736979166666666
$paramSheet.Range("E1").ColumnWidth = 16.451822916666668
$paramSheet.Range("F1").ColumnWidth = 22.877604166666668

# Select B19 on the PARAMETERS sheet, matching the commit's saved view state
$paramSheet.Range("B19").Select()

# Make PARAMETERS the active/selected tab
$paramSheet.Activate()

$wb.Save()
